$d = $word.ActiveDocument

function Replace-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex).Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $p.InsertXML($pkg)
}

# 1) "export PYTHONPATH..." paragraph: bold paragraph mark + bold/yellow-highlight run
$xml1 = '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:highlight w:val="yellow"/></w:rPr>' +
    '<w:t>export PYTHONPATH=$PYTHONPATH:`pwd`:`pwd`/slim</w:t></w:r></w:p>'
Replace-ParaXml 22 $xml1

# 2) protoc paragraph: add _GoBack bookmark at the end
$xml2 = '<w:p><w:r><w:t>/home/julyedu_433249/work/tf_base/research/bin/protoc object_detection/protos/*.proto --python_out=.</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Replace-ParaXml 24 $xml2

# 3) Merge "Tf." + "TFRecordReader" runs into a single run
$xml3 = '<w:p><w:pPr><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console" w:cs="Lucida Console"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console" w:cs="Lucida Console"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' +
    '<w:t>Tf.TFRecordReader</w:t></w:r></w:p>'
Replace-ParaXml 58 $xml3

# 4) "原尺寸:" paragraph: remove _GoBack bookmark (it moved earlier, to the protoc paragraph)
$xml4 = '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>原尺寸</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r></w:p>'
Replace-ParaXml 65 $xml4

# 5) "6.8 summary" paragraph: drop the paragraph-mark rPr (rFonts hint=eastAsia)
$xml5 = '<w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>6</w:t></w:r>' +
    '<w:r><w:t>.8 summary</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>实验</w:t></w:r></w:p>'
Replace-ParaXml 79 $xml5

Write-Output "done"
